$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.266.04'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '2.026.39'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.71'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.05'
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.378'
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.102'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '2.324.05'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.29'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.21'
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.739'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.19'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '2.029.70'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").Value = '37.145.81'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.73'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").Value = '0.0₃0817'
$ws.Range("E21").Value = '  -1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.00'
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.46'
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  -1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.93'
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.131'
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.70'
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.46'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0606'
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.00'
$ws.Range("E34").Value = '  +9.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.43'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.21'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.468.97'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("B41").Value = 'FTXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.36'
$ws.Range("E41").Value = '  +17.24%  '
$ws.Range("E42").Value = '  -2.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.81'
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.98'
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0915'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.20'
$ws.Range("E46").Value = '  -4.81%  '
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.11'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.93'
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").Value = '2.210.15'
$ws.Range("E51").Value = '  +0.14%  '
